$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 123, shifting the existing data (old rows
# 123-187) down to rows 124-188, and updating the sheet dimension
# automatically.
$ws.Rows(123).Insert()

# Populate the newly inserted row 123 with the new weekly record (same
# market/category template as the surrounding rows).
$ws.Range("A123").Value = 8
$ws.Range("B123").Value = "Terminal La Palmera de La Serena"
$ws.Range("C123").Value = "Coquimbo"
$ws.Range("D123").Value = 44529
$ws.Range("E123").Value = 4
$ws.Range("F123").Value = 100112012
$ws.Range("G123").Value = "Espinaca"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 2300
$ws.Range("K123").Value = 400
$ws.Range("L123").Value = 500
$ws.Range("M123").Value = 450
$ws.Range("N123").Value = "$/atado 300 a 500 gramos"
$ws.Range("O123").Value = "Provincia del Elquí"
$ws.Range("P123").Value = 900
$ws.Range("Q123").Value = 0.5
$ws.Range("R123").Value = "Hortaliza"
